$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: left-rotate C,D,E -> new C=old D, new D=old E, new E=old C
$c3 = $ws.Range("C3").Value2
$d3 = $ws.Range("D3").Value2
$e3 = $ws.Range("E3").Value2
$ws.Range("C3").Value2 = $d3
$ws.Range("D3").Value2 = $e3
$ws.Range("E3").Value2 = $c3

# Row 4: right-rotate C,D,E -> new C=old E, new D=old C, new E=old D
$c4 = $ws.Range("C4").Value2
$d4 = $ws.Range("D4").Value2
$e4 = $ws.Range("E4").Value2
$ws.Range("C4").Value2 = $e4
$ws.Range("D4").Value2 = $c4
$ws.Range("E4").Value2 = $d4

# Row 23: right-rotate C,D,E -> new C=old E, new D=old C, new E=old D
$c23 = $ws.Range("C23").Value2
$d23 = $ws.Range("D23").Value2
$e23 = $ws.Range("E23").Value2
$ws.Range("C23").Value2 = $e23
$ws.Range("D23").Value2 = $c23
$ws.Range("E23").Value2 = $d23

# Row 24: swap D,E
$d24 = $ws.Range("D24").Value2
$e24 = $ws.Range("E24").Value2
$ws.Range("D24").Value2 = $e24
$ws.Range("E24").Value2 = $d24

# Row 25: swap D,E
$d25 = $ws.Range("D25").Value2
$e25 = $ws.Range("E25").Value2
$ws.Range("D25").Value2 = $e25
$ws.Range("E25").Value2 = $d25

# Row 40: swap C,D
$c40 = $ws.Range("C40").Value2
$d40 = $ws.Range("D40").Value2
$ws.Range("C40").Value2 = $d40
$ws.Range("D40").Value2 = $c40

# Row 41: swap C,D
$c41 = $ws.Range("C41").Value2
$d41 = $ws.Range("D41").Value2
$ws.Range("C41").Value2 = $d41
$ws.Range("D41").Value2 = $c41

# Row 55: swap C,D
$c55 = $ws.Range("C55").Value2
$d55 = $ws.Range("D55").Value2
$ws.Range("C55").Value2 = $d55
$ws.Range("D55").Value2 = $c55
